$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header cell for 08-jul, mirroring the formatting of the previous
# day's header cell (X1) which carries the bold/border/centered style.
$ws.Range("Y1").Value = "08-jul"
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New day's hourly prices (column Y, rows 2-25)
$ws.Range("Y2").Value = 61.85
$ws.Range("Y3").Value = 43.45
$ws.Range("Y4").Value = 39.17
$ws.Range("Y5").Value = 31.81
$ws.Range("Y6").Value = 31.91
$ws.Range("Y7").Value = 30.1
$ws.Range("Y8").Value = 38.52
$ws.Range("Y9").Value = 61.13
$ws.Range("Y10").Value = 60
$ws.Range("Y11").Value = 18.01
$ws.Range("Y12").Value = 4.31
$ws.Range("Y13").Value = 1.72
$ws.Range("Y14").Value = 1.72
$ws.Range("Y15").Value = 0
$ws.Range("Y16").Value = 0
$ws.Range("Y17").Value = 0
$ws.Range("Y18").Value = 0
$ws.Range("Y19").Value = 11.78
$ws.Range("Y20").Value = 45
$ws.Range("Y21").Value = 60.39
$ws.Range("Y22").Value = 86.87
$ws.Range("Y23").Value = 82.46
$ws.Range("Y24").Value = 101.59
$ws.Range("Y25").Value = 88.16
